$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-04-22 Tuesday"; new = "2025-04-23 Wednesday"},
    @{old = "905×6="; new = "698×3="},
    @{old = "351×6="; new = "209×4="},
    @{old = "966×7="; new = "114×4="},
    @{old = "714×8="; new = "437×2="},
    @{old = "326×9="; new = "369×6="},
    @{old = "513×9="; new = "543×8="},
    @{old = "156×9="; new = "204×6="},
    @{old = "593×2="; new = "178×2="},
    @{old = "687×7="; new = "910×2="},
    @{old = "887×2="; new = "400×3="},
    @{old = "602×2="; new = "349×3="},
    @{old = "941×4="; new = "206×2="},
    @{old = "798×3="; new = "454×6="},
    @{old = "816×7="; new = "923×9="},
    @{old = "563×5="; new = "353×2="},
    @{old = "364×4="; new = "315×6="},
    @{old = "172×5="; new = "450×4="},
    @{old = "679×3="; new = "418×3="},
    @{old = "399×7="; new = "254×5="},
    @{old = "983×4="; new = "981×7="},
    @{old = "739×9="; new = "109×8="},
    @{old = "629×8="; new = "230×6="},
    @{old = "599×9="; new = "996×9="},
    @{old = "925×9="; new = "678×9="},
    @{old = "379×6="; new = "690×3="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
